# Hc-C5ar1.xlsx NATMI TPM refresh: drop stale rows 4-7 (old ECs / MuSCs
# combinations) and update rows 2-3 with newly recomputed TPM-derived
# specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four trailing result rows (old dimension A1:T7 -> A1:T3).
$ws.Rows("4:7").Delete()

# Row 2 (FAPs -> Hc -> C5ar1 -> ...): target cluster + recomputed metrics.
$ws.Range("D2").Value = "FAPs"
$ws.Range("I2").Value = 0.8736649195182647
$ws.Range("J2").Value = 0.8736649195182647
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1949616666666667
$ws.Range("N2").Value = 0.584885
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.05799628677888889
$ws.Range("R2").Value = 0.5219665810099999
$ws.Range("S2").Value = 0.8736649195182647
$ws.Range("T2").Value = 0.8736649195182647

# Row 3 (MuSCs -> Hc -> C5ar1 -> FAPs): sending cluster + recomputed metrics.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("G3").Value = 0.043016
$ws.Range("H3").Value = 0.129048
$ws.Range("I3").Value = 0.1263350804817352
$ws.Range("J3").Value = 0.1263350804817352
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.008386471053333334
$ws.Range("R3").Value = 0.07547823948
$ws.Range("S3").Value = 0.1263350804817352
$ws.Range("T3").Value = 0.1263350804817352
